# Update "想去人数" (number of people interested) figures that changed
# between crawl runs, as captured in the upstream data diff.
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)      -- no changes in this sheet
# Sheet 3 = 本地生活 (Local Life)
# Sheet 4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsLocalLife  = $wb.Worksheets.Item(3)   # 本地生活
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (sheet1) ---
$wsExhibition.Range("F3").Value  = 8663   # was 8473
$wsExhibition.Range("F8").Value  = 586    # was 585
$wsExhibition.Range("F16").Value = 8656   # was 8654
$wsExhibition.Range("F33").Value = 2149   # was 2147
$wsExhibition.Range("F39").Value = 219    # was 217
$wsExhibition.Range("F40").Value = 158    # was 157

# --- 本地生活 (sheet3) ---
$wsLocalLife.Range("F3").Value = 713      # was 712

# --- 全部类型 (sheet4) ---
$wsAllTypes.Range("F3").Value  = 713      # was 712
$wsAllTypes.Range("F6").Value  = 8663     # was 8473
$wsAllTypes.Range("F13").Value = 586      # was 585
$wsAllTypes.Range("F20").Value = 8656     # was 8654
$wsAllTypes.Range("F34").Value = 2149     # was 2147
$wsAllTypes.Range("F40").Value = 219      # was 217
$wsAllTypes.Range("F41").Value = 158      # was 157
